$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H2").Value = 162
$ws.Range("I2").Value = 104.55556
$ws.Range("K2").Value = 104.55556
$ws.Range("M2").Value = 8.44444
$ws.Range("H33").Value = 13814734
$ws.Range("I33").Value = 17267380
$ws.Range("K33").Value = 17267380
$ws.Range("M33").Value = -17267151
$ws.Range("H53").Value = 629.4
$ws.Range("J53").Value = 239
$ws.Range("L53").Value = 239
$ws.Range("N53").Value = -1513
$ws.Range("H58").Value = 1308.5454
$ws.Range("I58").Value = 432.33334
$ws.Range("K58").Value = 1297.00002
$ws.Range("M58").Value = -1147.00002
$ws.Range("H80").Value = 527494.75
$ws.Range("I80").Value = 562.8
$ws.Range("J80").Value = 715684.7
$ws.Range("K80").Value = 1688.4
$ws.Range("L80").Value = 2147054.1
$ws.Range("M80").Value = -690.3999999999999
$ws.Range("N80").Value = -2149050.1
$ws.Range("H83").Value = 527494.75
$ws.Range("I83").Value = 562.8
$ws.Range("J83").Value = 715684.7
$ws.Range("K83").Value = 5065.2
$ws.Range("L83").Value = 6441162.3
$ws.Range("M83").Value = -73.19999999999982
$ws.Range("N83").Value = -6451146.3
$ws.Range("H96").Value = 683.75
$ws.Range("I96").Value = 510.2143
$ws.Range("J96").Value = 1898.5
$ws.Range("K96").Value = 1530.6429
$ws.Range("L96").Value = 5695.5
$ws.Range("M96").Value = -157.6428999999998
$ws.Range("N96").Value = -8441.5
$ws.Range("H100").Value = 31368.572
$ws.Range("I100").Value = 43575.625
$ws.Range("J100").Value = 4735
$ws.Range("K100").Value = 43575.625
$ws.Range("L100").Value = 4735
$ws.Range("M100").Value = -43034.625
$ws.Range("N100").Value = -5817
$ws.Range("H101").Value = 596.25
$ws.Range("I101").Value = 642.8570999999999
$ws.Range("J101").Value = 270
$ws.Range("K101").Value = 1928.5713
$ws.Range("L101").Value = 810
$ws.Range("M101").Value = -306.5712999999998
$ws.Range("N101").Value = -4054
$ws.Range("H113").Value = 5167
$ws.Range("I113").Value = 2334
$ws.Range("K113").Value = 2334
$ws.Range("M113").Value = 920
$ws.Range("H129").Value = 1414.4231
$ws.Range("I129").Value = 988.7646999999999
$ws.Range("K129").Value = 2966.2941
$ws.Range("M129").Value = 2033.7059
$ws.Range("H132").Value = 369.66666
$ws.Range("I132").Value = 369.66666
$ws.Range("K132").Value = 1108.99998
$ws.Range("M132").Value = 1421.00002

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 5141.25
$ws.Range("I2").Value = 5133.6
$ws.Range("J2").Value = 5148.9
$ws.Range("K2").Value = 5133.6
$ws.Range("L2").Value = 5148.9
$ws.Range("M2").Value = -5020.6
$ws.Range("N2").Value = -5374.9
$ws.Range("H11").Value = 10075.143
$ws.Range("I11").Value = 18342.334
$ws.Range("J11").Value = 3874.75
$ws.Range("K11").Value = 18342.334
$ws.Range("L11").Value = 3874.75
$ws.Range("M11").Value = -18198.334
$ws.Range("N11").Value = -4162.75
$ws.Range("H35").Value = 15100.375
$ws.Range("I35").Value = 3953
$ws.Range("K35").Value = 3953
$ws.Range("M35").Value = -3547
$ws.Range("H45").Value = 7767.864
$ws.Range("I45").Value = 8433.157999999999
$ws.Range("K45").Value = 8433.157999999999
$ws.Range("M45").Value = -8056.157999999999
$ws.Range("H61").Value = 6155.421
$ws.Range("I61").Value = 5050.2856
$ws.Range("K61").Value = 5050.2856
$ws.Range("M61").Value = -4838.2856
$ws.Range("H74").Value = 12771.667
$ws.Range("I74").Value = 13725.84
$ws.Range("J74").Value = 8000.8
$ws.Range("K74").Value = 13725.84
$ws.Range("L74").Value = 8000.8
$ws.Range("M74").Value = -12851.84
$ws.Range("N74").Value = -9748.799999999999
$ws.Range("H75").Value = 68249.5
$ws.Range("J75").Value = 89332.664
$ws.Range("L75").Value = 89332.664
$ws.Range("N75").Value = -91080.664
$ws.Range("H77").Value = 12771.667
$ws.Range("I77").Value = 13725.84
$ws.Range("J77").Value = 8000.8
$ws.Range("K77").Value = 68629.2
$ws.Range("L77").Value = 40004
$ws.Range("M77").Value = -64261.2
$ws.Range("N77").Value = -48740
$ws.Range("H78").Value = 68249.5
$ws.Range("J78").Value = 89332.664
$ws.Range("L78").Value = 267997.992
$ws.Range("N78").Value = -276733.992
$ws.Range("H110").Value = 1250.409
$ws.Range("I110").Value = 1195.1765
$ws.Range("J110").Value = 1438.2
$ws.Range("K110").Value = 1195.1765
$ws.Range("L110").Value = 1438.2
$ws.Range("M110").Value = 849.8235
$ws.Range("N110").Value = -5528.2
$ws.Range("H116").Value = 5141.25
$ws.Range("I116").Value = 5133.6
$ws.Range("J116").Value = 5148.9
$ws.Range("K116").Value = 5133.6
$ws.Range("L116").Value = 5148.9
$ws.Range("M116").Value = -2839.6
$ws.Range("N116").Value = -9736.9
$ws.Range("H122").Value = 1160.2858
$ws.Range("I122").Value = 844.4
$ws.Range("J122").Value = 1950
$ws.Range("K122").Value = 2533.2
$ws.Range("L122").Value = 5850
$ws.Range("M122").Value = -83.19999999999982
$ws.Range("N122").Value = -10750
$ws.Range("H124").Value = 32499.5
$ws.Range("J124").Value = 32499.5
$ws.Range("L124").Value = 32499.5
$ws.Range("N124").Value = -42319.5
$ws.Range("H132").Value = 866.4286
$ws.Range("I132").Value = 847.04346
$ws.Range("K132").Value = 2541.13038
$ws.Range("M132").Value = -11.13038000000006
$ws.Range("H136").Value = 6155.421
$ws.Range("I136").Value = 5050.2856
$ws.Range("K136").Value = 15150.8568
$ws.Range("M136").Value = -12600.8568

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 5141.25
$ws.Range("I3").Value = 5133.6
$ws.Range("J3").Value = 5148.9
$ws.Range("K3").Value = 5133.6
$ws.Range("L3").Value = 5148.9
$ws.Range("M3").Value = -5019.6
$ws.Range("N3").Value = -5376.9
$ws.Range("H64").Value = 253.875
$ws.Range("J64").Value = 266.42856
$ws.Range("L64").Value = 266.42856
$ws.Range("N64").Value = -716.4285600000001
$ws.Range("H67").Value = 253.875
$ws.Range("J67").Value = 266.42856
$ws.Range("L67").Value = 266.42856
$ws.Range("N67").Value = -1826.42856
$ws.Range("H94").Value = 1631.8572
$ws.Range("I94").Value = 988.2308
$ws.Range("K94").Value = 988.2308
$ws.Range("M94").Value = -537.2308
$ws.Range("H105").Value = 3695.3547
$ws.Range("I105").Value = 3324.4075
$ws.Range("K105").Value = 3324.4075
$ws.Range("M105").Value = -1577.4075
$ws.Range("H130").Value = 34999
$ws.Range("J130").Value = 34999
$ws.Range("L130").Value = 34999
$ws.Range("N130").Value = -45039
$ws.Range("H131").Value = 40013
$ws.Range("J131").Value = 40013
$ws.Range("L131").Value = 40013
$ws.Range("N131").Value = -50093
$ws.Range("H132").Value = 70311.56
$ws.Range("J132").Value = 70311.56
$ws.Range("L132").Value = 70311.56
$ws.Range("N132").Value = -80431.56

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H22").Value = 230
$ws.Range("I22").Value = 229.04762
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 229.04762
$ws.Range("L22").Value = 250
$ws.Range("M22").Value = 120.95238
$ws.Range("N22").Value = -950
$ws.Range("H31").Value = 2525.926
$ws.Range("I31").Value = 1585.0526
$ws.Range("K31").Value = 1585.0526
$ws.Range("M31").Value = -1290.0526
$ws.Range("H34").Value = 2525.926
$ws.Range("I34").Value = 1585.0526
$ws.Range("K34").Value = 1585.0526
$ws.Range("M34").Value = -1383.0526
$ws.Range("H94").Value = 8485.571
$ws.Range("I94").Value = 7866.3335
$ws.Range("K94").Value = 7866.3335
$ws.Range("M94").Value = -7415.3335
$ws.Range("H132").Value = 23395.963
$ws.Range("I132").Value = 15021.29
$ws.Range("J132").Value = 34683.566
$ws.Range("K132").Value = 45063.87
$ws.Range("L132").Value = 104050.698
$ws.Range("M132").Value = -42533.87
$ws.Range("N132").Value = -109110.698
$ws.Range("H134").Value = 2840.2246
$ws.Range("I134").Value = 1717.7
$ws.Range("J134").Value = 7829.222
$ws.Range("K134").Value = 5153.1
$ws.Range("L134").Value = 23487.666
$ws.Range("M134").Value = -2618.1
$ws.Range("N134").Value = -28557.666

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H2").Value = 60.933334
$ws.Range("I2").Value = 43.857143
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 263.142858
$ws.Range("L2").Value = 1800
$ws.Range("M2").Value = -150.142858
$ws.Range("N2").Value = -2026
$ws.Range("H5").Value = 1607.9615
$ws.Range("I5").Value = 683.6667
$ws.Range("J5").Value = 1885.25
$ws.Range("K5").Value = 2051.0001
$ws.Range("L5").Value = 5655.75
$ws.Range("M5").Value = -1939.0001
$ws.Range("N5").Value = -5879.75
$ws.Range("H12").Value = 1697.6666
$ws.Range("J12").Value = 2032.1
$ws.Range("L12").Value = 6096.299999999999
$ws.Range("N12").Value = -6442.299999999999
$ws.Range("H40").Value = 48.0625
$ws.Range("J40").Value = 50.333332
$ws.Range("L40").Value = 201.333328
$ws.Range("N40").Value = -339.333328
$ws.Range("H55").Value = 4845.636
$ws.Range("I55").Value = 2832.3333
$ws.Range("J55").Value = 5600.625
$ws.Range("K55").Value = 8496.999899999999
$ws.Range("L55").Value = 16801.875
$ws.Range("M55").Value = -8319.999899999999
$ws.Range("N55").Value = -17155.875
$ws.Range("H87").Value = 5640
$ws.Range("I87").Value = 5640
$ws.Range("K87").Value = 16920
$ws.Range("M87").Value = -15672
$ws.Range("H90").Value = 5640
$ws.Range("I90").Value = 5640
$ws.Range("K90").Value = 50760
$ws.Range("M90").Value = -44520
$ws.Range("H125").Value = 34676.668
$ws.Range("I125").Value = 100030
$ws.Range("J125").Value = 2000
$ws.Range("K125").Value = 300090
$ws.Range("L125").Value = 6000
$ws.Range("M125").Value = -295170
$ws.Range("N125").Value = -15840
$ws.Range("H129").Value = 2798.6
$ws.Range("I129").Value = 1330
$ws.Range("J129").Value = 3165.75
$ws.Range("K129").Value = 3990
$ws.Range("L129").Value = 9497.25
$ws.Range("M129").Value = 1010
$ws.Range("N129").Value = -19497.25
$ws.Range("H131").Value = 2275.9023
$ws.Range("I131").Value = 848.4666999999999
$ws.Range("J131").Value = 3099.423
$ws.Range("K131").Value = 2545.4001
$ws.Range("L131").Value = 9298.269
$ws.Range("M131").Value = 2494.5999
$ws.Range("N131").Value = -19378.269
$ws.Range("H133").Value = 5937.5
$ws.Range("J133").Value = 6250
$ws.Range("L133").Value = 18750
$ws.Range("N133").Value = -28870
$ws.Range("H135").Value = 1607.9615
$ws.Range("I135").Value = 683.6667
$ws.Range("J135").Value = 1885.25
$ws.Range("K135").Value = 6153.0003
$ws.Range("L135").Value = 16967.25
$ws.Range("M135").Value = -3618.0003
$ws.Range("N135").Value = -22037.25
$ws.Range("H140").Value = 8999.666999999999
$ws.Range("I140").Value = 8799.6
$ws.Range("K140").Value = 26398.8
$ws.Range("M140").Value = -21218.8

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H9").Value = 527.75
$ws.Range("I9").Value = 57
$ws.Range("J9").Value = 684.6667
$ws.Range("K9").Value = 57
$ws.Range("L9").Value = 684.6667
$ws.Range("M9").Value = 113
$ws.Range("N9").Value = -1024.6667
$ws.Range("H13").Value = 2919.8
$ws.Range("I13").Value = 1750
$ws.Range("J13").Value = 3699.6667
$ws.Range("K13").Value = 1750
$ws.Range("L13").Value = 3699.6667
$ws.Range("M13").Value = -1611
$ws.Range("N13").Value = -3977.6667
$ws.Range("H25").Value = 2500
$ws.Range("J25").Value = 2500
$ws.Range("L25").Value = 2500
$ws.Range("N25").Value = -3558
$ws.Range("H63").Value = 89998.5
$ws.Range("I63").Value = 89998
$ws.Range("K63").Value = 89998
$ws.Range("M63").Value = -89312
$ws.Range("H66").Value = 89998.5
$ws.Range("I66").Value = 89998
$ws.Range("K66").Value = 269994
$ws.Range("M66").Value = -266562
$ws.Range("H93").Value = 42994.6
$ws.Range("J93").Value = 42994.6
$ws.Range("L93").Value = 42994.6
$ws.Range("N93").Value = -46738.6
$ws.Range("H98").Value = 26634.5
$ws.Range("J98").Value = 26634.5
$ws.Range("L98").Value = 26634.5
$ws.Range("N98").Value = -32624.5
$ws.Range("H102").Value = 1694.1875
$ws.Range("I102").Value = 1722.2858
$ws.Range("K102").Value = 1722.2858
$ws.Range("M102").Value = -100.2858000000001
$ws.Range("H104").Value = 54832.668
$ws.Range("J104").Value = 54832.668
$ws.Range("L104").Value = 54832.668
$ws.Range("N104").Value = -61820.668
$ws.Range("H132").Value = 15576.275
$ws.Range("I132").Value = 14484.091
$ws.Range("J132").Value = 19008.857
$ws.Range("K132").Value = 43452.273
$ws.Range("L132").Value = 57026.571
$ws.Range("M132").Value = -40922.273
$ws.Range("N132").Value = -62086.571

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 1410.0889
$ws.Range("I16").Value = 1451.6857
$ws.Range("K16").Value = 1451.6857
$ws.Range("M16").Value = -1281.6857
$ws.Range("H22").Value = 1730
$ws.Range("J22").Value = 1796.9231
$ws.Range("L22").Value = 1796.9231
$ws.Range("N22").Value = -2386.9231
$ws.Range("H27").Value = 1730
$ws.Range("J27").Value = 1796.9231
$ws.Range("L27").Value = 1796.9231
$ws.Range("N27").Value = -2010.9231
$ws.Range("H46").Value = 1227.2041
$ws.Range("J46").Value = 1470.7188
$ws.Range("L46").Value = 1470.7188
$ws.Range("N46").Value = -1846.7188
$ws.Range("H55").Value = 415.6
$ws.Range("I55").Value = 253.8
$ws.Range("J55").Value = 577.4
$ws.Range("K55").Value = 253.8
$ws.Range("L55").Value = 577.4
$ws.Range("M55").Value = -80.80000000000001
$ws.Range("N55").Value = -923.4
$ws.Range("H61").Value = 5110.75
$ws.Range("I61").Value = 5285.2856
$ws.Range("J61").Value = 3889
$ws.Range("K61").Value = 5285.2856
$ws.Range("L61").Value = 3889
$ws.Range("M61").Value = -5083.2856
$ws.Range("N61").Value = -4293
$ws.Range("H63").Value = 43208.168
$ws.Range("J63").Value = 39850
$ws.Range("L63").Value = 39850
$ws.Range("N63").Value = -41348
$ws.Range("H66").Value = 43208.168
$ws.Range("J66").Value = 39850
$ws.Range("L66").Value = 119550
$ws.Range("N66").Value = -127038
$ws.Range("H76").Value = 18072
$ws.Range("J76").Value = 18072
$ws.Range("L76").Value = 18072
$ws.Range("N76").Value = -18748
$ws.Range("H79").Value = 18072
$ws.Range("J79").Value = 18072
$ws.Range("L79").Value = 18072
$ws.Range("N79").Value = -20412
$ws.Range("H93").Value = 3068.762
$ws.Range("J93").Value = 1155.6666
$ws.Range("L93").Value = 1155.6666
$ws.Range("N93").Value = -3651.6666
$ws.Range("H113").Value = 5110.75
$ws.Range("I113").Value = 5285.2856
$ws.Range("J113").Value = 3889
$ws.Range("K113").Value = 5285.2856
$ws.Range("L113").Value = 3889
$ws.Range("M113").Value = -3115.2856
$ws.Range("N113").Value = -8229
$ws.Range("H122").Value = 4687
$ws.Range("I122").Value = 5222.375
$ws.Range("K122").Value = 15667.125
$ws.Range("M122").Value = -13217.125
$ws.Range("H132").Value = 4800.796
$ws.Range("I132").Value = 4832.8945
$ws.Range("K132").Value = 14498.6835
$ws.Range("M132").Value = -11968.6835

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H23").Value = 49.666668
$ws.Range("I23").Value = 49.666668
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 49.666668
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 179.333332
$ws.Range("N23").ClearContents()
$ws.Range("H95").Value = 35000
$ws.Range("J95").Value = 35000
$ws.Range("L95").Value = 35000
$ws.Range("N95").Value = -40492
$ws.Range("H100").Value = 814.7083
$ws.Range("I100").Value = 867.7619
$ws.Range("J100").Value = 443.33334
$ws.Range("K100").Value = 1735.5238
$ws.Range("L100").Value = 886.66668
$ws.Range("M100").Value = -1194.5238
$ws.Range("N100").Value = -1968.66668
$ws.Range("H107").Value = 6174231
$ws.Range("I107").Value = 1036
$ws.Range("K107").Value = 3108
$ws.Range("M107").Value = -1188
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H124").Value = 32012.666
$ws.Range("J124").Value = 32012.666
$ws.Range("L124").Value = 32012.666
$ws.Range("N124").Value = -41832.666
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840
$ws.Range("H126").Value = 4065.415
$ws.Range("I126").Value = 2982.2974
$ws.Range("J126").Value = 6570.125
$ws.Range("K126").Value = 8946.8922
$ws.Range("L126").Value = 19710.375
$ws.Range("M126").Value = -6476.8922
$ws.Range("N126").Value = -24650.375
